# Pedidos.xlsx update
# - Insert 2 new order rows right after the header row (new rows 2 and 3),
#   pushing the existing order rows down by two.
# - Append 7 new order rows at the bottom of the existing data block
#   (rows 23-29 in the new layout).
# - Move the active selection to I13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as genuine text (shared string), even when the text
# looks like a plain number (e.g. "80265160"), without leaving behind a
# throw-away number-format/quote-prefix style. We do this by writing a
# formula that evaluates to the literal string, then collapsing the formula
# down to its value in place via copy / paste-special-values. The resulting
# cell keeps whatever style it already had and stores a plain text value.
function Set-TextValue {
    param($Row, $Col, [string]$Text)

    $cell = $ws.Cells.Item($Row, $Col)
    $escaped = $Text.Replace('"', '""')
    $cell.FormulaR1C1 = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

$ws.Application.CutCopyMode = $false

# --- Insert two new rows right after the header (row 1) ---------------------
$ws.Rows.Item(2).Resize(2).EntireRow.Insert()

# The freshly inserted rows pick up the header row's formatting; restore the
# plain data-row formatting (matching every other order row) by copying it
# from the row immediately below (the first untouched original data row).
$ws.Range("A4:C4").Copy()
$ws.Range("A2:C3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Row 2
Set-TextValue 2 1 "80265160"
Set-TextValue 2 2 "30018-KMT-I"
$ws.Cells.Item(2, 3).Value = 24000

# Row 3
Set-TextValue 3 1 "80265942"
Set-TextValue 3 2 "20637-TDK-I"
$ws.Cells.Item(3, 3).Value = 2000

# --- Append seven new rows at the end of the existing data (rows 23-29) -----
Set-TextValue 23 1 "80266510"
Set-TextValue 23 2 "10025-ARI-I"
$ws.Cells.Item(23, 3).Value = 1

Set-TextValue 24 1 "80266511"
Set-TextValue 24 2 "10020-ARI-I"
$ws.Cells.Item(24, 3).Value = 1

Set-TextValue 25 1 "80266512"
Set-TextValue 25 2 "10382-ARI-I"
$ws.Cells.Item(25, 3).Value = 1

Set-TextValue 26 1 "80266513"
Set-TextValue 26 2 "10399-ARI-I"
$ws.Cells.Item(26, 3).Value = 1

Set-TextValue 27 1 "80266514"
Set-TextValue 27 2 "10008-LDG-I"
$ws.Cells.Item(27, 3).Value = 1

Set-TextValue 28 1 "80266516"
Set-TextValue 28 2 "10001-LDG-I"
$ws.Cells.Item(28, 3).Value = 1

Set-TextValue 29 1 "84004823"
Set-TextValue 29 2 "15211-DLO-I"
$ws.Cells.Item(29, 3).Value = 1

# --- Update the active selection -------------------------------------------
$ws.Range("I13").Select()
